$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted right after the existing
# row for 2023-03-23 (row 44). Insert a new row at 45, which shifts the
# former rows 45-55 down to 46-56, then populate the new row 45 with the
# latest data point.
$ws.Rows.Item(45).Insert()

$ws.Range("A45").Value = 1
$ws.Range("B45").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C45").Value = "Arica y Parinacota"
$ws.Range("D45").Value = 45204
$ws.Range("E45").Value = 15
$ws.Range("F45").Value = 100112044
$ws.Range("G45").Value = "Perejil"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 200
$ws.Range("K45").Value = 1800
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = 1900
$ws.Range("N45").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O45").Value = "Región de Arica y Parinacota"
$ws.Range("P45").Value = 950
$ws.Range("Q45").Value = 2
$ws.Range("R45").Value = "Hortaliza"
